$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $vals) {
    $n = $vals.Count
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0,$i] = $vals[$i]
    }
    $rng = $ws.Range($ws.Cells.Item($r,1), $ws.Cells.Item($r,$n))
    $rng.Value2 = $arr
}

Set-Row 20 @(11, 'Vega Monumental Concepción', 'Bíobío', 44953, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Especial', 200, 12000, 12000, 12000, '$/bandeja 18 kilos granel', 'Provincia de Curicó', 667, 18)
Set-Row 21 @(11, 'Vega Monumental Concepción', 'Bíobío', 44953, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 200, 11000, 11000, 11000, '$/bandeja 18 kilos granel', 'Provincia de Curicó', 611, 18)
Set-Row 22 @(11, 'Vega Monumental Concepción', 'Bíobío', 44953, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 200, 10000, 10000, 10000, '$/bandeja 18 kilos granel', 'Provincia de Curicó', 556, 18)
Set-Row 23 @(11, 'Vega Monumental Concepción', 'Bíobío', 44236, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Lemon', 'Primera', 200, 14000, 15000, 14500, '$/caja 16 kilos granel', 'Región de O''Higgins', 906, 16)
Set-Row 24 @(11, 'Vega Monumental Concepción', 'Bíobío', 44236, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Lemon', 'Segunda', 100, 12000, 12000, 12000, '$/caja 16 kilos granel', 'Región de O''Higgins', 750, 16)
Set-Row 25 @(11, 'Vega Monumental Concepción', 'Bíobío', 44323, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 200, 11000, 12000, 11500, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 639, 18)
Set-Row 26 @(11, 'Vega Monumental Concepción', 'Bíobío', 44323, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Segunda', 100, 9000, 9000, 9000, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 500, 18)
Set-Row 27 @(11, 'Vega Monumental Concepción', 'Bíobío', 44622, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 270, 8000, 8500, 8222, '$/caja 16 kilos granel', 'Región de O''Higgins', 514, 16)
Set-Row 28 @(11, 'Vega Monumental Concepción', 'Bíobío', 44596, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 220, 7000, 8000, 7545, '$/caja 16 kilos granel', 'Provincia de Curicó', 472, 16)
Set-Row 29 @(11, 'Vega Monumental Concepción', 'Bíobío', 44952, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 150, 10000, 10000, 10000, '$/bandeja 18 kilos granel', 'Provincia de Curicó', 556, 18)
Set-Row 30 @(11, 'Vega Monumental Concepción', 'Bíobío', 44952, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 200, 8000, 8000, 8000, '$/bandeja 18 kilos granel', 'Provincia de Curicó', 444, 18)
Set-Row 31 @(11, 'Vega Monumental Concepción', 'Bíobío', 44952, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Fortuna', 'Primera', 270, 9000, 10000, 9556, '$/bandeja 18 kilos granel', 'Provincia de Curicó', 531, 18)
Set-Row 32 @(11, 'Vega Monumental Concepción', 'Bíobío', 44335, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 100, 10000, 11000, 10500, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 583, 18)
Set-Row 33 @(11, 'Vega Monumental Concepción', 'Bíobío', 44335, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Segunda', 50, 9000, 9000, 9000, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 500, 18)
Set-Row 34 @(11, 'Vega Monumental Concepción', 'Bíobío', 44615, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 200, 10000, 11000, 10500, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 583, 18)
Set-Row 35 @(11, 'Vega Monumental Concepción', 'Bíobío', 44615, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 100, 9000, 9000, 9000, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 500, 18)
Set-Row 36 @(11, 'Vega Monumental Concepción', 'Bíobío', 44285, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 200, 9000, 10000, 9500, '$/caja 18 kilos granel', 'Región de O''Higgins', 528, 18)
Set-Row 37 @(11, 'Vega Monumental Concepción', 'Bíobío', 44285, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Segunda', 100, 8000, 8000, 8000, '$/caja 18 kilos granel', 'Región de O''Higgins', 444, 18)
Set-Row 38 @(11, 'Vega Monumental Concepción', 'Bíobío', 44594, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 170, 9000, 9500, 9265, '$/caja 16 kilos granel', 'Provincia de Curicó', 579, 16)
Set-Row 39 @(11, 'Vega Monumental Concepción', 'Bíobío', 44589, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 220, 7500, 8000, 7773, '$/caja 16 kilos granel', 'Región de O''Higgins', 486, 16)
Set-Row 40 @(11, 'Vega Monumental Concepción', 'Bíobío', 44218, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 200, 10000, 11000, 10500, '$/caja 16 kilos granel', 'Región de O''Higgins', 656, 16)
Set-Row 41 @(11, 'Vega Monumental Concepción', 'Bíobío', 44218, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 100, 9000, 9000, 9000, '$/caja 16 kilos granel', 'Región de O''Higgins', 562, 16)
Set-Row 42 @(11, 'Vega Monumental Concepción', 'Bíobío', 44939, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Especial', 150, 14000, 14000, 14000, '$/caja 15 kilos empedrada', 'Región de O''Higgins', 933, 15)
Set-Row 43 @(11, 'Vega Monumental Concepción', 'Bíobío', 44939, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 120, 12000, 12000, 12000, '$/caja 15 kilos empedrada', 'Región de O''Higgins', 800, 15)
Set-Row 44 @(11, 'Vega Monumental Concepción', 'Bíobío', 44939, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 150, 11000, 11000, 11000, '$/caja 15 kilos empedrada', 'Región de O''Higgins', 733, 15)
Set-Row 45 @(11, 'Vega Monumental Concepción', 'Bíobío', 44939, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Larry Ann', 'Especial', 120, 14000, 14000, 14000, '$/caja 15 kilos empedrada', 'Región de O''Higgins', 933, 15)
Set-Row 46 @(11, 'Vega Monumental Concepción', 'Bíobío', 44939, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Larry Ann', 'Primera', 100, 12000, 12000, 12000, '$/caja 15 kilos empedrada', 'Región de O''Higgins', 800, 15)
Set-Row 47 @(11, 'Vega Monumental Concepción', 'Bíobío', 44939, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Larry Ann', 'Segunda', 120, 11000, 11000, 11000, '$/caja 15 kilos empedrada', 'Región de O''Higgins', 733, 15)
Set-Row 48 @(11, 'Vega Monumental Concepción', 'Bíobío', 44246, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 100, 10000, 10000, 10000, '$/caja 16 kilos granel', 'Región de O''Higgins', 625, 16)
Set-Row 49 @(11, 'Vega Monumental Concepción', 'Bíobío', 44246, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Segunda', 100, 8000, 8000, 8000, '$/caja 16 kilos granel', 'Región de O''Higgins', 500, 16)
Set-Row 50 @(11, 'Vega Monumental Concepción', 'Bíobío', 44580, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 230, 10000, 11000, 10565, '$/bandeja 18 kilos granel', 'Provincia de Curicó', 587, 18)
Set-Row 51 @(11, 'Vega Monumental Concepción', 'Bíobío', 44922, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 100, 13000, 14000, 13500, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 750, 18)
Set-Row 52 @(11, 'Vega Monumental Concepción', 'Bíobío', 44558, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 100, 17000, 18000, 17500, '$/bandeja 18 kilos granel', 'Provincia de Curicó', 972, 18)
Set-Row 53 @(11, 'Vega Monumental Concepción', 'Bíobío', 44558, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 50, 16000, 16000, 16000, '$/bandeja 18 kilos granel', 'Provincia de Curicó', 889, 18)
Set-Row 54 @(11, 'Vega Monumental Concepción', 'Bíobío', 44299, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 200, 12000, 13000, 12500, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 694, 18)
Set-Row 55 @(11, 'Vega Monumental Concepción', 'Bíobío', 44299, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Segunda', 100, 11000, 11000, 11000, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 611, 18)
Set-Row 56 @(11, 'Vega Monumental Concepción', 'Bíobío', 44642, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 220, 8000, 9000, 8455, '$/bandeja 18 kilos granel', 'Provincia de Curicó', 470, 18)
Set-Row 57 @(11, 'Vega Monumental Concepción', 'Bíobío', 44588, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 220, 9000, 9500, 9227, '$/caja 16 kilos granel', 'Región de O''Higgins', 577, 16)
Set-Row 58 @(11, 'Vega Monumental Concepción', 'Bíobío', 44202, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 200, 14000, 15000, 14500, '$/caja 18 kilos granel', 'Región de O''Higgins', 806, 18)
Set-Row 59 @(11, 'Vega Monumental Concepción', 'Bíobío', 44202, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 100, 12000, 12000, 12000, '$/caja 18 kilos granel', 'Región de O''Higgins', 667, 18)
Set-Row 60 @(11, 'Vega Monumental Concepción', 'Bíobío', 44931, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 100, 15000, 16000, 15500, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 861, 18)
Set-Row 61 @(11, 'Vega Monumental Concepción', 'Bíobío', 44931, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 50, 14000, 14000, 14000, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 778, 18)
Set-Row 62 @(11, 'Vega Monumental Concepción', 'Bíobío', 44343, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 200, 10000, 11000, 10500, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 583, 18)
Set-Row 63 @(11, 'Vega Monumental Concepción', 'Bíobío', 44343, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Segunda', 100, 9000, 9000, 9000, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 500, 18)
Set-Row 64 @(11, 'Vega Monumental Concepción', 'Bíobío', 44645, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 200, 8000, 8500, 8250, '$/bandeja 18 kilos granel', 'Provincia de Curicó', 458, 18)
Set-Row 65 @(11, 'Vega Monumental Concepción', 'Bíobío', 44645, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Segunda', 200, 6000, 6500, 6250, '$/bandeja 18 kilos granel', 'Provincia de Curicó', 347, 18)
Set-Row 66 @(11, 'Vega Monumental Concepción', 'Bíobío', 44586, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 200, 9000, 10000, 9500, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 528, 18)
Set-Row 67 @(11, 'Vega Monumental Concepción', 'Bíobío', 44586, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 100, 8000, 8000, 8000, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 444, 18)
Set-Row 68 @(11, 'Vega Monumental Concepción', 'Bíobío', 44307, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 200, 9000, 10000, 9500, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 528, 18)
Set-Row 69 @(11, 'Vega Monumental Concepción', 'Bíobío', 44307, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Segunda', 100, 8000, 8000, 8000, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 444, 18)
Set-Row 70 @(11, 'Vega Monumental Concepción', 'Bíobío', 44215, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 200, 10000, 11000, 10500, '$/caja 16 kilos granel', 'Región de O''Higgins', 656, 16)
Set-Row 71 @(11, 'Vega Monumental Concepción', 'Bíobío', 44215, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 100, 8000, 8000, 8000, '$/caja 16 kilos granel', 'Región de O''Higgins', 500, 16)
Set-Row 72 @(11, 'Vega Monumental Concepción', 'Bíobío', 44951, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Fortuna', 'Primera', 170, 10000, 11000, 10529, '$/caja 18 kilos granel', 'Provincia de Curicó', 585, 18)
Set-Row 73 @(11, 'Vega Monumental Concepción', 'Bíobío', 44266, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 200, 9000, 10000, 9500, '$/caja 18 kilos granel', 'Región de O''Higgins', 528, 18)
Set-Row 74 @(11, 'Vega Monumental Concepción', 'Bíobío', 44266, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 100, 8000, 8000, 8000, '$/caja 18 kilos granel', 'Región de O''Higgins', 444, 18)
Set-Row 75 @(11, 'Vega Monumental Concepción', 'Bíobío', 44607, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 250, 11000, 12000, 11520, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 640, 18)
Set-Row 76 @(11, 'Vega Monumental Concepción', 'Bíobío', 44607, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 170, 9000, 9500, 9265, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 515, 18)
Set-Row 77 @(11, 'Vega Monumental Concepción', 'Bíobío', 44328, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 100, 9000, 10000, 9500, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 528, 18)
Set-Row 78 @(11, 'Vega Monumental Concepción', 'Bíobío', 44328, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Segunda', 50, 8000, 8000, 8000, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 444, 18)
Set-Row 79 @(11, 'Vega Monumental Concepción', 'Bíobío', 44946, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 100, 11000, 12000, 11500, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 639, 18)
Set-Row 80 @(11, 'Vega Monumental Concepción', 'Bíobío', 44946, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 50, 9000, 9000, 9000, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 500, 18)
Set-Row 81 @(11, 'Vega Monumental Concepción', 'Bíobío', 44644, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 350, 8500, 9000, 8714, '$/bandeja 18 kilos granel', 'Provincia de Curicó', 484, 18)
Set-Row 82 @(11, 'Vega Monumental Concepción', 'Bíobío', 44637, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 220, 8000, 8500, 8273, '$/bandeja 18 kilos granel', 'Provincia de Curicó', 460, 18)
Set-Row 83 @(11, 'Vega Monumental Concepción', 'Bíobío', 44223, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 100, 10000, 11000, 10500, '$/caja 16 kilos granel', 'Región de O''Higgins', 656, 16)
Set-Row 84 @(11, 'Vega Monumental Concepción', 'Bíobío', 44223, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 50, 9000, 9000, 9000, '$/caja 16 kilos granel', 'Región de O''Higgins', 562, 16)
Set-Row 85 @(11, 'Vega Monumental Concepción', 'Bíobío', 44616, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 200, 9000, 10000, 9500, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 528, 18)
Set-Row 86 @(11, 'Vega Monumental Concepción', 'Bíobío', 44616, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 100, 8000, 8000, 8000, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 444, 18)
Set-Row 87 @(11, 'Vega Monumental Concepción', 'Bíobío', 44631, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 100, 8000, 9000, 8500, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 472, 18)
Set-Row 88 @(11, 'Vega Monumental Concepción', 'Bíobío', 44631, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Segunda', 50, 7000, 7000, 7000, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 389, 18)
Set-Row 89 @(11, 'Vega Monumental Concepción', 'Bíobío', 44643, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 180, 8500, 9000, 8778, '$/bandeja 18 kilos granel', 'Región del Maule', 488, 18)
Set-Row 90 @(11, 'Vega Monumental Concepción', 'Bíobío', 44643, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Segunda', 150, 6500, 7000, 6767, '$/bandeja 18 kilos granel', 'Región del Maule', 376, 18)
Set-Row 91 @(11, 'Vega Monumental Concepción', 'Bíobío', 44657, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 200, 9000, 10000, 9500, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 528, 18)
Set-Row 92 @(11, 'Vega Monumental Concepción', 'Bíobío', 44657, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Segunda', 100, 8000, 8000, 8000, '$/bandeja 18 kilos granel', 'Región de O''Higgins', 444, 18)
Set-Row 93 @(11, 'Vega Monumental Concepción', 'Bíobío', 44251, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 200, 9000, 10000, 9500, '$/caja 16 kilos granel', 'Región de O''Higgins', 594, 16)
Set-Row 94 @(11, 'Vega Monumental Concepción', 'Bíobío', 44251, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Segunda', 100, 8000, 8000, 8000, '$/caja 16 kilos granel', 'Región de O''Higgins', 500, 16)
Set-Row 95 @(11, 'Vega Monumental Concepción', 'Bíobío', 44636, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Segunda', 220, 6500, 7000, 6727, '$/bandeja 18 kilos granel', 'Provincia de Curicó', 374, 18)
Set-Row 96 @(11, 'Vega Monumental Concepción', 'Bíobío', 44595, 8, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 250, 8500, 9000, 8740, '$/caja 16 kilos granel', 'Provincia de Curicó', 546, 16)

# Rows 94-96 are brand new rows; give column D (Fecha) the same date/time
# number format used by the rest of the date column (style index 2 in the
# original workbook -> numFmtId 165 "YYYY-MM-DD HH:MM:SS").
$ws.Range("D94:D96").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Refresh the sheet's used-range dimension to A1:T96.
$ws.Range("A1:T96").Select()
